$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 162
$ws.Range("F5").Value = 46
$ws.Range("F6").Value = 2731
$ws.Range("F8").Value = 1614
$ws.Range("F9").Value = 7402
$ws.Range("F11").Value = 7582
$ws.Range("F12").Value = 17
$ws.Range("F14").Value = 5
$ws.Range("F15").Value = 6068
$ws.Range("F16").Value = 3239
$ws.Range("F17").Value = 3606
$ws.Range("F18").Value = 12
$ws.Range("F19").Value = 4
$ws.Range("F21").Value = 26
$ws.Range("F22").Value = 439
$ws.Range("F23").Value = 3
$ws.Range("F25").Value = 278
$ws.Range("F26").Value = 2096
$ws.Range("F27").Value = 112
$ws.Range("F29").Value = 925
$ws.Range("F30").Value = 254
$ws.Range("F31").Value = 1064
$ws.Range("F33").Value = 12
$ws.Range("F34").Value = 2588
$ws.Range("F35").Value = 1430
$ws.Range("F36").Value = 5
$ws.Range("F38").Value = 14
$ws.Range("F39").Value = 3191
$ws.Range("F40").Value = 144
$ws.Range("F41").Value = 238
$ws.Range("F43").Value = 891
$ws.Range("F45").Value = 1250

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 9
$ws.Range("F6").Value = 43
$ws.Range("F9").Value = 395
$ws.Range("F18").Value = 11

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9
$ws.Range("F5").Value = 162
$ws.Range("F7").Value = 46
$ws.Range("F9").Value = 2731
$ws.Range("F10").Value = 1614
$ws.Range("F12").Value = 43
$ws.Range("F13").Value = 7402
$ws.Range("F14").Value = 7582
$ws.Range("F15").Value = 17
$ws.Range("F16").Value = 5
$ws.Range("F17").Value = 6068
$ws.Range("F18").Value = 3239
$ws.Range("F19").Value = 3606
$ws.Range("F20").Value = 12
$ws.Range("F21").Value = 4
$ws.Range("F22").Value = 26
$ws.Range("F23").Value = 439
$ws.Range("F24").Value = 3
$ws.Range("F28").Value = 278
$ws.Range("F29").Value = 2096
$ws.Range("F31").Value = 112
$ws.Range("F35").Value = 925
$ws.Range("F36").Value = 254
$ws.Range("F37").Value = 12
$ws.Range("F38").Value = 2588
$ws.Range("F39").Value = 1430
$ws.Range("F40").Value = 5
$ws.Range("F43").Value = 3191
$ws.Range("F44").Value = 238
$ws.Range("F45").Value = 891
$ws.Range("F47").Value = 1250
